# Implements a Dreadnought (large, single-model unit) into the 40k sim workbook,
# and fixes an inverse-AP bug on the existing melee weapons.

$wb = $excel.ActiveWorkbook

$wsRanged = $wb.Worksheets.Item("Templar Ranged Weapons")
$wsMelee  = $wb.Worksheets.Item("Templar Melee Weapons")
$wsModels = $wb.Worksheets.Item("Templar Models")

# --- Templar Ranged Weapons: add the Dreadnought's Assault Cannon ---
$rangedRow = New-Object 'object[,]' 1,9
$rangedRow[0,0] = "Assault Cannon"  # Weapon
$rangedRow[0,1] = 24                # Range
$rangedRow[0,2] = "Heavy"           # Weapon Type
$rangedRow[0,3] = 0                 # Die/Fixed
$rangedRow[0,4] = 6                 # Shots
$rangedRow[0,5] = 6                 # S
$rangedRow[0,6] = 1                 # AP
$rangedRow[0,7] = 0                 # Die/Fixed
$rangedRow[0,8] = 1                 # D
$wsRanged.Range("A30:I30").Value = $rangedRow

# --- Templar Models: add the Dreadnought unit ---
$modelRow = New-Object 'object[,]' 1,10
$modelRow[0,0] = "Dreadnought"  # UNIT NAME
$modelRow[0,1] = 6              # M
$modelRow[0,2] = 3              # WS
$modelRow[0,3] = 3              # BS
$modelRow[0,4] = 6              # S
$modelRow[0,5] = 7              # T
$modelRow[0,6] = 8              # W
$modelRow[0,7] = 4              # A
$modelRow[0,8] = 8              # Ld
$modelRow[0,9] = 3              # Sv
$wsModels.Range("A17:J17").Value = $modelRow
$wsModels.Range("L17").Value = 30  # Radius

# --- Templar Melee Weapons: fix inverse AP bug (AP should be positive here) ---
$wsMelee.Range("C5").Value = 3
$wsMelee.Range("C6").Value = 1

# --- Templar Melee Weapons: add the Dreadnought Combat Weapon ---
$meleeRow = New-Object 'object[,]' 1,5
$meleeRow[0,0] = "Dreadnought Combat Weapon"  # Weapon
$meleeRow[0,1] = "2x"                          # S
$meleeRow[0,2] = 3                             # AP
$meleeRow[0,3] = 0                             # Die/Fixed
$meleeRow[0,4] = 3                             # D
$wsMelee.Range("A8:E8").Value = $meleeRow

# --- Restore view/selection state on each sheet, ending with the melee sheet active ---
$wsRanged.Range("A31").Select()
$wsModels.Range("J8").Select()
$wsMelee.Range("C8").Select()
